# Update res_bus vm_pu results for the 380 kV case (Case_4_195)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-bus voltage magnitude values for data rows 2-25 (Excel bus index 0-23),
# one array per column letter, in row order.
$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$colNumbers = @{
    "B" = 2
    "C" = 3
    "D" = 4
    "E" = 5
    "F" = 6
    "I" = 9
    "J" = 10
    "K" = 11
    "L" = 12
    "M" = 13
    "N" = 14
}

$values = @{
    "B" = @(1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02)
    "C" = @(1.052438541728832,1.053605160018208,1.054359975593028,1.054677285756479,1.054730562723162,1.05436421556176,1.052832819871069,1.050133760179284,1.048333946441267,1.047554485279862,1.047264937566294,1.04732704751597,1.047530551615969,1.047655934359897,1.048385673836425,1.048843384207577,1.049110346494656,1.049201371676467,1.048794277513405,1.047470625327055,1.046638269593525,1.047079529321466,1.048816466752654,1.050831604004078)
    "D" = @(1.05811653983289,1.059051079870336,1.059655353741973,1.059909286271834,1.059951916555361,1.059658747210119,1.05843246240599,1.056268243529698,1.054823163598723,1.054196885562752,1.053964175021989,1.054014095987131,1.054177651332981,1.054278412166208,1.054864715999961,1.05523234158136,1.055446718439305,1.055519806447001,1.055192904313448,1.054129490655443,1.053460399841218,1.053815143260542,1.05521072448052,1.056828144407369)
    "E" = @(1.065480819714477,1.066604246093982,1.067331404611604,1.067637156829645,1.067688497200394,1.067335489871797,1.065860440345208,1.063262924382915,1.061532358080504,1.06078325766702,1.060505044467161,1.060564720526948,1.060760259747255,1.060880742662791,1.061582078514915,1.062022073221036,1.062278738508895,1.062366258808102,1.061974863495931,1.060702677333086,1.059903011487016,1.060326910074835,1.061996195464642,1.063934245849413)
    "F" = @(1.071575003673295,1.072705262244345,1.073436734134581,1.073744273912814,1.073795912909669,1.073440843381602,1.071956955510614,1.069343043700957,1.067600996900605,1.066846795358711,1.066566667878478,1.066626755439477,1.066823639601399,1.066944948597766,1.067651053144347,1.06809400434856,1.068352381566501,1.068440483460619,1.068046478716535,1.066765661673292,1.065960456480527,1.066387302258114,1.068067953468747,1.070018700991869)
    "I" = @(1.04318107110154,1.043449227149097,1.04362084525048,1.043692539397632,1.043704550531321,1.043621805015519,1.043272088434486,1.042641315672917,1.042211030590215,1.042022393910516,1.041951977042902,1.041967097487369,1.042016580346931,1.042047022143878,1.042223500880967,1.042333579959686,1.04239756335217,1.042419342068762,1.042321792673707,1.042002018508429,1.041798945617991,1.041906789765167,1.042327119532456,1.042806107362481)
    "J" = @(1.057461142390235,1.058277188405152,1.058804608385093,1.059026188640041,1.059063384326659,1.058807569728772,1.057737056857172,1.055845934898941,1.054581966962958,1.054033883159373,1.053830182849989,1.053873882591434,1.054017047631488,1.054105240733206,1.054618325034506,1.054939960585253,1.055127490211269,1.055191420278937,1.054905459892471,1.053974892347708,1.053389126912124,1.053699717040386,1.054921049501757,1.056335400425)
    "K" = @(1.060849653076864,1.061597857029376,1.062080932490133,1.062283763298383,1.062317804589456,1.06208364372352,1.061102732265271,1.059366090991897,1.058202830325939,1.057697815659244,1.057510032555787,1.057550321646076,1.057682297495378,1.057763585894655,1.058236318727987,1.058532499036789,1.05870512914961,1.058763970058688,1.058500734817469,1.057643439372405,1.057103276850494,1.057389736130122,1.058515088092585,1.059816022254732)
    "L" = @(1.068193944902512,1.069132008682469,1.06973863480127,1.069993573789368,1.070036374100971,1.069742041648327,1.068511043581803,1.066339058687364,1.064889146041955,1.06426085327521,1.064027405753488,1.064077484293684,1.064241557896054,1.064342639542662,1.06493083376753,1.065299665762632,1.065514753835509,1.065588085636789,1.06526009828094,1.0641932442933,1.063522056090927,1.063877905255194,1.065277977262682,1.066900905214148)
    "M" = @(1.074271812736264,1.075217876375621,1.075829571415243,1.076086615586815,1.076129767854143,1.075833006493268,1.074591637239294,1.072400555638437,1.070937362456106,1.070303189335274,1.070067538131518,1.070118090248758,1.070283712195736,1.070385745251641,1.070979437710525,1.071351683883883,1.071568750863069,1.071642755294286,1.071311751388492,1.070234943180842,1.06955738381574,1.069916621108016,1.071329795342918,1.072967436161269)
    "N" = @(1.058962858504917,1.059780063398778,1.060308232375585,1.060530127299912,1.060567375808671,1.060311197924713,1.059239164802045,1.057345357233719,1.056079594318288,1.055530732172836,1.055326742585647,1.055370504385738,1.05551387273657,1.055602191082608,1.056116004022468,1.056438096332589,1.056625892272183,1.056689913127881,1.056403546644865,1.055471657587452,1.054885060297761,1.055196091499631,1.056419158393185,1.057835517856968)
}

foreach ($col in $columns) {
    $colNum = $colNumbers[$col]
    $colValues = $values[$col]
    for ($i = 0; $i -lt $colValues.Length; $i++) {
        $row = 2 + $i
        $ws.Cells.Item($row, $colNum).Value = $colValues[$i]
    }
}
